$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 17
$ws.Cells.Item(17, 8).Value = 302304.12
$ws.Cells.Item(17, 10).Value = 302304.12
$ws.Cells.Item(17, 12).Value = 906912.36
$ws.Cells.Item(17, 14).Value = -907248.36

# ALC row 58
$ws.Cells.Item(58, 8).Value = 3362.5
$ws.Cells.Item(58, 9).Value = 225
$ws.Cells.Item(58, 11).Value = 675
$ws.Cells.Item(58, 13).Value = -525

# ALC row 62
$ws.Cells.Item(62, 8).Value = 97536.73
$ws.Cells.Item(62, 9).Value = 147429.14
$ws.Cells.Item(62, 10).Value = 10225
$ws.Cells.Item(62, 11).Value = 147429.14
$ws.Cells.Item(62, 12).Value = 10225
$ws.Cells.Item(62, 13).Value = -146805.14
$ws.Cells.Item(62, 14).Value = -11473

# ALC row 65
$ws.Cells.Item(65, 8).Value = 97536.73
$ws.Cells.Item(65, 9).Value = 147429.14
$ws.Cells.Item(65, 10).Value = 10225
$ws.Cells.Item(65, 11).Value = 737145.7000000001
$ws.Cells.Item(65, 12).Value = 51125
$ws.Cells.Item(65, 13).Value = -734025.7000000001
$ws.Cells.Item(65, 14).Value = -57365

# ALC row 106
$ws.Cells.Item(106, 8).Value = 3021.1538
$ws.Cells.Item(106, 9).Value = 3196.875
$ws.Cells.Item(106, 11).Value = 3196.875
$ws.Cells.Item(106, 13).Value = -2565.875

# ALC row 107
$ws.Cells.Item(107, 8).Value = 205.71428
$ws.Cells.Item(107, 9).Value = 106.666664
$ws.Cells.Item(107, 11).Value = 106.666664
$ws.Cells.Item(107, 13).Value = 1813.333336

# ALC row 123
$ws.Cells.Item(123, 8).Value = 20695.422
$ws.Cells.Item(123, 10).Value = 20695.422
$ws.Cells.Item(123, 12).Value = 20695.422
$ws.Cells.Item(123, 14).Value = -30495.422

# ALC row 132
$ws.Cells.Item(132, 8).Value = 1842.2295
$ws.Cells.Item(132, 9).Value = 1720.08
$ws.Cells.Item(132, 10).Value = 2397.4546
$ws.Cells.Item(132, 11).Value = 5160.24
$ws.Cells.Item(132, 12).Value = 7192.3638
$ws.Cells.Item(132, 13).Value = -2630.24
$ws.Cells.Item(132, 14).Value = -12252.3638

# ALC row 138
$ws.Cells.Item(138, 8).Value = 3042.65
$ws.Cells.Item(138, 9).Value = 722.7406999999999
$ws.Cells.Item(138, 10).Value = 3900.6987
$ws.Cells.Item(138, 11).Value = 2168.2221
$ws.Cells.Item(138, 12).Value = 11702.0961
$ws.Cells.Item(138, 13).Value = 2971.7779
$ws.Cells.Item(138, 14).Value = -21982.0961

# ALC row 141
$ws.Cells.Item(141, 8).Value = 3681.111
$ws.Cells.Item(141, 9).Value = 2278.182
$ws.Cells.Item(141, 10).Value = 5885.7144
$ws.Cells.Item(141, 11).Value = 6834.545999999999
$ws.Cells.Item(141, 12).Value = 17657.1432
$ws.Cells.Item(141, 13).Value = -1654.545999999999
$ws.Cells.Item(141, 14).Value = -28017.1432

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Cells.Item(32, 8).Value = 8553.333000000001
$ws.Cells.Item(32, 9).Value = 4490
$ws.Cells.Item(32, 10).Value = 44542.855
$ws.Cells.Item(32, 11).Value = 4490
$ws.Cells.Item(32, 12).Value = 44542.855
$ws.Cells.Item(32, 13).Value = -4203
$ws.Cells.Item(32, 14).Value = -45116.855

# ARM row 63
$ws.Cells.Item(63, 8).Value = 1975
$ws.Cells.Item(63, 9).Value = 1828.5714
$ws.Cells.Item(63, 10).Value = 3000
$ws.Cells.Item(63, 11).Value = 1828.5714
$ws.Cells.Item(63, 12).Value = 3000
$ws.Cells.Item(63, 13).Value = -1142.5714
$ws.Cells.Item(63, 14).Value = -4372

# ARM row 66
$ws.Cells.Item(66, 8).Value = 1975
$ws.Cells.Item(66, 9).Value = 1828.5714
$ws.Cells.Item(66, 10).Value = 3000
$ws.Cells.Item(66, 11).Value = 9142.857
$ws.Cells.Item(66, 12).Value = 15000
$ws.Cells.Item(66, 13).Value = -5710.857
$ws.Cells.Item(66, 14).Value = -21864

# ARM row 125
$ws.Cells.Item(125, 8).Value = 33122.5
$ws.Cells.Item(125, 10).Value = 33122.5
$ws.Cells.Item(125, 12).Value = 33122.5
$ws.Cells.Item(125, 14).Value = -42962.5

$ws = $wb.Worksheets.Item("BSM")
# BSM row 26
$ws.Cells.Item(26, 8).Value = 9900
$ws.Cells.Item(26, 9).Value = 9900
$ws.Cells.Item(26, 11).Value = 9900
$ws.Cells.Item(26, 13).Value = -9608

# BSM row 64
$ws.Cells.Item(64, 8).Value = 906.7778
$ws.Cells.Item(64, 9).Value = 383.33334
$ws.Cells.Item(64, 10).Value = 1168.5
$ws.Cells.Item(64, 11).Value = 383.33334
$ws.Cells.Item(64, 12).Value = 1168.5
$ws.Cells.Item(64, 13).Value = -158.33334
$ws.Cells.Item(64, 14).Value = -1618.5

# BSM row 67
$ws.Cells.Item(67, 8).Value = 906.7778
$ws.Cells.Item(67, 9).Value = 383.33334
$ws.Cells.Item(67, 10).Value = 1168.5
$ws.Cells.Item(67, 11).Value = 383.33334
$ws.Cells.Item(67, 12).Value = 1168.5
$ws.Cells.Item(67, 13).Value = 396.66666
$ws.Cells.Item(67, 14).Value = -2728.5

$ws = $wb.Worksheets.Item("CRP")
# CRP row 99
$ws.Cells.Item(99, 8).Value = 2362.6775
$ws.Cells.Item(99, 9).Value = 1966.45
$ws.Cells.Item(99, 10).Value = 3083.0908
$ws.Cells.Item(99, 11).Value = 1966.45
$ws.Cells.Item(99, 12).Value = 3083.0908
$ws.Cells.Item(99, 13).Value = -468.45
$ws.Cells.Item(99, 14).Value = -6079.0908

# CRP row 105
$ws.Cells.Item(105, 8).Value = 1456.5264
$ws.Cells.Item(105, 9).Value = 1273.375
$ws.Cells.Item(105, 10).Value = 2433.3333
$ws.Cells.Item(105, 11).Value = 1273.375
$ws.Cells.Item(105, 12).Value = 2433.3333
$ws.Cells.Item(105, 13).Value = 473.625
$ws.Cells.Item(105, 14).Value = -5927.3333

# CRP row 126
$ws.Cells.Item(126, 8).Value = 2362.6775
$ws.Cells.Item(126, 9).Value = 1966.45
$ws.Cells.Item(126, 10).Value = 3083.0908
$ws.Cells.Item(126, 11).Value = 5899.35
$ws.Cells.Item(126, 12).Value = 9249.2724
$ws.Cells.Item(126, 13).Value = -3429.35
$ws.Cells.Item(126, 14).Value = -14189.2724

# CRP row 134
$ws.Cells.Item(134, 8).Value = 2704.2646
$ws.Cells.Item(134, 9).Value = 1947.75
$ws.Cells.Item(134, 10).Value = 4519.9
$ws.Cells.Item(134, 11).Value = 5843.25
$ws.Cells.Item(134, 12).Value = 13559.7
$ws.Cells.Item(134, 13).Value = -3308.25
$ws.Cells.Item(134, 14).Value = -18629.7

$ws = $wb.Worksheets.Item("CUL")
# CUL row 96
$ws.Cells.Item(96, 8).Value = 5000
$ws.Cells.Item(96, 10).Value = 5000
$ws.Cells.Item(96, 12).Value = 15000
$ws.Cells.Item(96, 14).Value = -19118

# CUL row 117
$ws.Cells.Item(117, 8).Value = 2142.318
$ws.Cells.Item(117, 9).Value = 472.5
$ws.Cells.Item(117, 10).Value = 2768.5
$ws.Cells.Item(117, 11).Value = 1417.5
$ws.Cells.Item(117, 12).Value = 8305.5
$ws.Cells.Item(117, 13).Value = 2024.5
$ws.Cells.Item(117, 14).Value = -15189.5

# CUL row 141
$ws.Cells.Item(141, 8).Value = 6641.8184
$ws.Cells.Item(141, 9).Value = 5351.8184
$ws.Cells.Item(141, 10).Value = 7931.8184
$ws.Cells.Item(141, 11).Value = 16055.4552
$ws.Cells.Item(141, 12).Value = 23795.4552
$ws.Cells.Item(141, 13).Value = -10875.4552
$ws.Cells.Item(141, 14).Value = -34155.4552

$ws = $wb.Worksheets.Item("GSM")
# GSM row 70
$ws.Cells.Item(70, 8).Value = 6227.25
$ws.Cells.Item(70, 9).Value = 4900
$ws.Cells.Item(70, 11).Value = 4900
$ws.Cells.Item(70, 13).Value = -4630

# GSM row 73
$ws.Cells.Item(73, 8).Value = 6227.25
$ws.Cells.Item(73, 9).Value = 4900
$ws.Cells.Item(73, 11).Value = 4900
$ws.Cells.Item(73, 13).Value = -3964

# GSM row 136
$ws.Cells.Item(136, 8).Value = 22163
$ws.Cells.Item(136, 10).Value = 22163
$ws.Cells.Item(136, 12).Value = 66489
$ws.Cells.Item(136, 14).Value = -71589

$ws = $wb.Worksheets.Item("LTW")
# LTW row 40
$ws.Cells.Item(40, 8).Value = 1263985.8
$ws.Cells.Item(40, 9).Value = 1263985.8
$ws.Cells.Item(40, 11).Value = 1263985.8
$ws.Cells.Item(40, 13).Value = -1263849.8

# LTW row 64
$ws.Cells.Item(64, 8).Value = 150
$ws.Cells.Item(64, 10).Value = 150
$ws.Cells.Item(64, 12).Value = 150
$ws.Cells.Item(64, 14).Value = -600

# LTW row 67
$ws.Cells.Item(67, 8).Value = 150
$ws.Cells.Item(67, 10).Value = 150
$ws.Cells.Item(67, 12).Value = 150
$ws.Cells.Item(67, 14).Value = -1710

# LTW row 122
$ws.Cells.Item(122, 8).Value = 4706.3237
$ws.Cells.Item(122, 9).Value = 5169.231
$ws.Cells.Item(122, 11).Value = 15507.693
$ws.Cells.Item(122, 13).Value = -13057.693

# LTW row 132
$ws.Cells.Item(132, 8).Value = 1872.9032
$ws.Cells.Item(132, 9).Value = 1383.9048
$ws.Cells.Item(132, 10).Value = 2899.8
$ws.Cells.Item(132, 11).Value = 4151.7144
$ws.Cells.Item(132, 12).Value = 8699.400000000001
$ws.Cells.Item(132, 13).Value = -1621.7144
$ws.Cells.Item(132, 14).Value = -13759.4

$ws = $wb.Worksheets.Item("WVR")
# WVR row 63
$ws.Cells.Item(63, 8).Value = 2000
$ws.Cells.Item(63, 10).Value = 2000
$ws.Cells.Item(63, 12).Value = 2000
$ws.Cells.Item(63, 14).Value = -3248

# WVR row 66
$ws.Cells.Item(66, 8).Value = 2000
$ws.Cells.Item(66, 10).Value = 2000
$ws.Cells.Item(66, 12).Value = 6000
$ws.Cells.Item(66, 14).Value = -12240

# WVR row 99
$ws.Cells.Item(99, 8).Value = 0
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 14).ClearContents()

# WVR row 126
$ws.Cells.Item(126, 8).Value = 688.0714
$ws.Cells.Item(126, 9).Value = 652.75
$ws.Cells.Item(126, 10).Value = 900
$ws.Cells.Item(126, 11).Value = 1958.25
$ws.Cells.Item(126, 12).Value = 2700
$ws.Cells.Item(126, 13).Value = 511.75
$ws.Cells.Item(126, 14).Value = -7640
